$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# This shared string is reused by the Overview summary row (E2/F2) and by the
# per-file "Status" column (C2) on both language sheets, so update every
# occurrence to keep them pointing at the same (new) shared string.
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: report was regenerated, handback datetime refreshed ---
$wsZhCn.Range("K2").Value = "2016-08-17 16:48:32"

# --- zh-cn / de-de sheets: stale "handback not latest" warning cleared ---
$wsZhCn.Range("P2").Value = "'"
$wsZhCn.Range("P2").Style = "Normal"
$wsDeDe.Range("P2").Value = "'"
$wsDeDe.Range("P2").Style = "Normal"

# --- de-de sheet: report was regenerated, handback datetime refreshed ---
$wsDeDe.Range("K2").Value = "2016-08-17 16:48:39"

# --- Column width adjustments (report column auto-sizing) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

$wsZhCn.Columns.Item(3).ColumnWidth = 29.14
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8

$wsDeDe.Columns.Item(3).ColumnWidth = 29.14
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8
